$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preserving formatting like
# trailing/leading zeros that Excel would otherwise "clean up" if the
# string were auto-detected as a number. We temporarily force a text
# number format, assign the value, then restore the cell's original
# style so no stray formatting is left behind.
function Set-TextValue($Range, $Value) {
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

# ---- Price (column D) updates ----
Set-TextValue $ws.Range("D2")  "245.11"
Set-TextValue $ws.Range("D4")  "5.216"
Set-TextValue $ws.Range("D5")  "0.05797"
Set-TextValue $ws.Range("D6")  "6.515"
Set-TextValue $ws.Range("D8")  "0.8158"
Set-TextValue $ws.Range("D9")  "0.8587"
Set-TextValue $ws.Range("D10") "0.1363"
Set-TextValue $ws.Range("D11") "0.06973"
Set-TextValue $ws.Range("D12") "0.03181"
Set-TextValue $ws.Range("D13") "0.02868"
Set-TextValue $ws.Range("D14") "0.09373"
Set-TextValue $ws.Range("D15") "3.756"
Set-TextValue $ws.Range("D16") "0.001523"
Set-TextValue $ws.Range("D17") "0.04705"

Set-TextValue $ws.Range("D18") "0.0005971"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws.Range("D19") "0.006274"
Set-TextValue $ws.Range("D20") "0.001234"
Set-TextValue $ws.Range("D21") "0.004535"
Set-TextValue $ws.Range("D22") "0.00006903"
Set-TextValue $ws.Range("D25") "0.3177"
Set-TextValue $ws.Range("D28") "0.0002329"
Set-TextValue $ws.Range("D40") "0.03656"

# ---- Rows 41-43: coin ranking shuffled ----
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006293"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1052"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003038"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

Set-TextValue $ws.Range("D44") "0.007951"
Set-TextValue $ws.Range("D45") "0.00005279"
Set-TextValue $ws.Range("D47") "0.3400"
Set-TextValue $ws.Range("D48") "0.002342"
Set-TextValue $ws.Range("D50") "0.0002001"
